$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to remain as plain text, matching the original
# inlineStr cell type (prevents Excel from auto-converting numeric-looking
# strings such as "1.007" or "26.32" into actual numbers).
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '30.001.70'
$ws.Range('E2').Value = '  +0.33%  '

$ws.Range('D3').Value = '1.906.08'
$ws.Range('E3').Value = '  +0.10%  '

$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.55%  '

$ws.Range('D5').Value = '0.7908'
$ws.Range('E5').Value = '  -1.08%  '

$ws.Range('D6').Value = '242.98'
$ws.Range('E6').Value = '  +0.77%  '

$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  +0.26%  '

$ws.Range('D8').Value = '0.3202'
$ws.Range('E8').Value = '  +2.56%  '

$ws.Range('D9').Value = '26.32'
$ws.Range('E9').Value = '  +0.52%  '

$ws.Range('D10').Value = '0.07102'
$ws.Range('E10').Value = '  +3.34%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.08058'
$ws.Range('E11').Value = '  +0.97%  '

$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = '0.7721'
$ws.Range('E12').Value = '  +5.15%  '

$ws.Range('D13').Value = '1.920.28'
$ws.Range('E13').Value = '  +0.18%  '

$ws.Range('D14').Value = '5.330'
$ws.Range('E14').Value = '  +3.20%  '

$ws.Range('D15').Value = '92.52'
$ws.Range('E15').Value = '  -0.15%  '

$ws.Range('D16').Value = '30.052.50'
$ws.Range('E16').Value = '  +0.48%  '

$ws.Range('D17').Value = '13.96'
$ws.Range('E17').Value = '  +0.49%  '

$ws.Range('D18').Value = '5.974'
$ws.Range('E18').Value = '  +2.13%  '

$ws.Range('D19').Value = '246.04'
$ws.Range('E19').Value = '  +0.44%  '

$ws.Range('D20').Value = '0.000007732'
$ws.Range('E20').Value = '  +0.61%  '

$ws.Range('D21').Value = '2.163.12'
$ws.Range('E21').Value = '  -0.19%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '8.123'
$ws.Range('E23').Value = '  +17.64%  '

$ws.Range('D24').Value = '1.009'
$ws.Range('E24').Value = '  +0.77%  '

$ws.Range('D25').Value = '0.1607'
$ws.Range('E25').Value = '  +13.80%  '

$ws.Range('D26').Value = '9.346'
$ws.Range('E26').Value = '  +1.83%  '

$ws.Range('D27').Value = '166.05'
$ws.Range('E27').Value = '  -1.18%  '

$ws.Range('D28').Value = '18.74'
$ws.Range('E28').Value = '  -0.55%  '

$ws.Range('D29').Value = '2.113'
$ws.Range('E29').Value = '  +4.97%  '

$ws.Range('D30').Value = '1.383'
$ws.Range('E30').Value = '  +1.80%  '

$ws.Range('D31').Value = '1.543'
$ws.Range('E31').Value = '  +1.68%  '

$ws.Range('D32').Value = '4.497'
$ws.Range('E32').Value = '  +5.07%  '

$ws.Range('D33').Value = '0.05691'
$ws.Range('E33').Value = '  +2.96%  '

$ws.Range('D34').Value = '4.085'
$ws.Range('E34').Value = '  +0.70%  '

$ws.Range('D35').Value = '1.271'
$ws.Range('E35').Value = '  +1.37%  '

$ws.Range('D36').Value = '0.7382'
$ws.Range('E36').Value = '  +1.34%  '

$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '2.734'
$ws.Range('E37').Value = '  +0.24%  '

$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').Value = '0.9973'
$ws.Range('E38').Value = '  -0.35%  '

$ws.Range('D39').Value = '0.01932'
$ws.Range('E39').Value = '  +0.53%  '

$ws.Range('D40').Value = '2.790'
$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('D41').Value = '0.4450'
$ws.Range('E41').Value = '  +1.34%  '

$ws.Range('D42').Value = '72.82'
$ws.Range('E42').Value = '  +1.26%  '

$ws.Range('D43').Value = '5.966'
$ws.Range('E43').Value = '  -2.62%  '

$ws.Range('D44').Value = '0.8468'
$ws.Range('E44').Value = '  +1.38%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '1.003'
$ws.Range('E45').Value = '  +0.13%  '

$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '1.892'
$ws.Range('E46').Value = '  +1.63%  '

$ws.Range('D47').Value = '1.029.19'
$ws.Range('E47').Value = '  +5.34%  '

$ws.Range('D48').Value = '102.05'
$ws.Range('E48').Value = '  +1.49%  '

$ws.Range('D49').Value = '9.949'
$ws.Range('E49').Value = '  +2.50%  '

$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').Value = '3.047'
$ws.Range('E50').Value = '  +10.28%  '

$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '7.504'
$ws.Range('E51').Value = '  -0.42%  '
